$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>' + `
  '<w:p>' + `
  '<w:r/>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' + `
  '<w:r><w:t>: Read our review of Caishen''s Fortune XL, an Asian-themed slot game with 28 paylines, free spins, and a maximum jackpot of 3000x. Play for free now!</w:t></w:r>' + `
  '</w:p>' + `
  '</w:body></w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'

$metaRange = $d.Range($metaPara.Range.Start, $metaPara.Range.End)
$metaRange.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# 2) Remove the duplicated bold title paragraph near the end of the document
#    (it duplicated the H1 text) and rewrite the italic "meta description"
#    paragraph that follows it into the new image-prompt paragraph, keeping
#    its italic formatting.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($count - 1)
$dupTitlePara.Range.Delete()

$imgPromptPara = $d.Paragraphs.Item($d.Paragraphs.Count)

$imgXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>' + `
  '<w:p>' + `
  '<w:r><w:rPr><w:i/></w:rPr><w:t>Create a feature image for "Caishen''s Fortune XL" that captures the vibrant and luxurious Chinese theme of the game. The image should be in a cartoon style and feature a happy Maya warrior wearing glasses, symbolizing the exciting adventure and good fortune that players can experience in this online slot game. The background of the image should showcase traditional Chinese designs and colors, with cherry blossoms and gazebos. The image should be visually stunning and eye-catching, drawing players to try their luck with Cai Shen''s Fortune XL.</w:t></w:r>' + `
  '</w:p>' + `
  '</w:body></w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'

$imgRange = $d.Range($imgPromptPara.Range.Start, $imgPromptPara.Range.End)
$imgRange.InsertXML($imgXml)

Write-Output "done"
